$wb = $excel.ActiveWorkbook

# New handoff attempt identifiers (the previous handoff transform failed, so
# a brand-new handoff file id is generated and the status is updated to
# reflect the failed transform; the stale "Ready for handoff" specific
# handoff-file / handoff-datetime / reason details are cleared back to their
# "never handed off" defaults).
$newFileName = "0e1993c0-23dc-4c91-9f34-364b5fb4460d.md"
$newStatus   = "Handoff transform failed"
$zeroDate    = "0001-01-01 00:00:00"
$ignored     = "Ignored"

$fileNameUrl = "https://github.com/OpenLocalizationTest/oltest/blob/85283f7defaa6cd31d91550410230fe350635384/e2e/$newFileName"
$configUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/85283f7defaa6cd31d91550410230fe350635384/.localization-config"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# Rebuild hyperlinks (engine's Hyperlinks.Delete affects the whole sheet, so
# remove them all and re-add with the updated display text / same targets).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $fileNameUrl, "", "", $newFileName) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de)
# ---------------------------------------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A2").Value = $newFileName
    $ws.Range("B2").Value = $newStatus

    # Handoff transform failed, so there is no specific handoff file anymore.
    $ws.Range("C2").Clear()

    # Reset handoff/handback bookkeeping back to "never happened" state and
    # flag the row as ignored (same as row 3's ignored dependency row).
    $ws.Range("D2").Value = $zeroDate
    $ws.Range("G2").Value = $zeroDate
    $ws.Range("H2").Value = $ignored

    $ws.Range("D3").Value = $zeroDate
    $ws.Range("G3").Value = $zeroDate
    $ws.Range("H3").Value = $ignored

    # Rebuild hyperlinks: drop the old C2 (handoff file) hyperlink entirely,
    # keep A2/A3 but refresh the A2 display text.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $fileNameUrl, "", "", $newFileName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null
}
